# Add a new ordinal outcome row ("O23a" / ordinal_v2a) to the derived
# variables table, right after the existing "O23" (ordinal_v2) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# The table's header is row 1, and the existing "O23" row is worksheet
# row 133, so the new row needs to land on worksheet row 134 (shifting
# everything from the old row 134 ("Rx01") onward down by one row).
$insertRow = 134

$ws.Rows.Item($insertRow).Insert() | Out-Null

$ws.Range("A$insertRow").Value = "O23a"
$ws.Range("B$insertRow").Value = "ordinal_v2a"
$ws.Range("C$insertRow").Value = "Outcome"
$ws.Range("D$insertRow").Value = "Custom ordinal including need for oxygen in the hospital"

# Grow the table / autofilter range to include the freshly inserted row.
$lastRow = $lo.Range.Rows.Count + 1
$lo.Resize($ws.Range("A1:E$lastRow")) | Out-Null

# Reflect the author's final cursor position/selection.
$ws.Range("A135").Select() | Out-Null
